# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bad Drivers section
$ws.Range("D3").Value = 93.2
$ws.Range("C4").Value = 8409
$ws.Range("D4").Value = 96.90000000000001
$ws.Range("C5").Value = 10214

# Good Drivers section
$ws.Range("B13").Value = 449371
$ws.Range("B17").Value = 77999
